# asset profile changes commit
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update values in row 2 of the AssetProfile sheet
# (order chosen to match shared-string insertion order of the target file)
$ws.Cells.Item(2, 7).Value = "ddd - ddd"         # G2
$ws.Cells.Item(2, 8).Value = "sss - sss"         # H2
$ws.Cells.Item(2, 6).Value = "aaa"               # F2
$ws.Cells.Item(2, 1).Value = "Auto-13"          # A2

# F2 gets its own style: text number format, left/top alignment, keep existing border
$f2 = $ws.Cells.Item(2, 6)
$f2.NumberFormat = "@"
$f2.HorizontalAlignment = -4131   # xlHAlignLeft
$f2.VerticalAlignment = -4160    # xlVAlignTop

# Update selection on the active sheet to A2
$ws.Range("A2").Select()
